$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(300, 240, 520, 670, 100, 1000, 2000, 2350, 4500, 4570)

$row = 9
foreach ($v in $values) {
    $ws.Cells.Item($row, 1).Value = $v
    $row++
}

$ws.Range("A18").Select()
